# CIERRE 28 AGOSTO 21
# Advance the payroll receipt sheet from "SEMANA 34 (16-22 AGOSTO 2021)" to
# "SEMANA 35 (23-29 AGOSTO 2021)" and update this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The week banner text lives once as a literal cell (B9); every other
# occurrence (B27, H9, H27, B43, H43, B60) is a formula that chases back to
# B9, so a single write here ripples through all of them.
$ws.Range("B9").Value = "SEMANA   35  DEL  23  Al    29  AGOSTO         2021"

# --- Block 1 (rows 3-6) --------------------------------------------------
# INFONAVIT discount for this employee goes to 0 this week.
$ws.Range("K4").Value = 0

# --- Block 2 (rows 21-24) -------------------------------------------------
# This employee worked 6 extra-hour units this week (previously none), and
# the "extra" label that used to flag J22 no longer applies.
$ws.Range("J21").Value = 6
$ws.Range("K21").Value = 1680
$ws.Range("J22").ClearContents()

# --- Block 5 (rows 38-41) -------------------------------------------------
# PRESTAMO payment recorded this week.
$ws.Range("K40").Value = 833

# --- Selection cursor (cosmetic) ------------------------------------------
$ws.Range("I60:I61").Select()

$wb.Save()
